# Validated power modelling for cyclists
# Update the header row of each power-curve data sheet to use descriptive
# text labels instead of bare numbers, then leave the workbook focused on
# the metadata sheet (matching the author's final view state).

$wb = $excel.ActiveWorkbook

$dataSheetNames = @("M W", "M W kg", "F W", "F W kg")

foreach ($name in $dataSheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B1").Value = "5s [W]"
    $ws.Range("C1").Value = "60s [W]"
    $ws.Range("D1").Value = "300s [W]"
    $ws.Range("E1").Value = "1200s [W]"
    $ws.Range("A1").Value = "percentile [%]"

    # Reset selection back to the top-left cell on each of these sheets.
    $ws.Select() | Out-Null
    $ws.Range("A1").Select() | Out-Null
}

# Return focus to the metadata sheet, selecting B1 (the hyperlink cell),
# matching the saved view state in the workbook.
$metaSheet = $wb.Worksheets.Item("metadata")
$metaSheet.Select() | Out-Null
$metaSheet.Range("B1").Select() | Out-Null
